# Improved mappings and example for Coverage writes to buffer file (#386)
#
# This script applies the following changes to Sheet1:
#  - Row 42 (60.01 PATIENT NAME): update the Notes (G) text with a link to the
#    R4PatientTransformer source (identifier section).
#  - Rows 49, 50, 54 (60.08 INSURED'S DOB / 60.09 INSURED'S SSN /
#    60.13 INSURED'S SEX): Contained Resource changes from
#    "Patient/RelatedPerson" to "RelatedPerson", and new Fhir Path / Notes
#    values are added.
#  - Rows 51, 52 (60.10 / 60.11): Contained Resource stays "Patient" (no
#    change in text, only shared-string index shifts internally).
#  - Rows 71-78 (62.02 - 62.09 SUBSCRIBER ADDRESS / PHONE): Contained
#    Resource changes from "Patient/RelatedPerson" to "RelatedPerson", and a
#    Fhir Path value (.address / .telecom) is added.
#  - Row 102 (91.01 NAME OF INSURED): Contained Resource changes from
#    "Patient/RelatedPerson" to "RelatedPerson", and a Fhir Path value
#    (.name.text) is added.
#  - Column G is widened to fit the long URLs now stored in Notes.
#  - The active selection/view is moved to reflect where the author was
#    working (F103).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: PATIENT NAME -------------------------------------------------
$ws.Range('G42').Value = 'See R4PatientTransformer for identifier system/code details: https://github.com/department-of-veterans-affairs/health-apis-data-query/blob/a6b5bb26a88606753fdfff227cd23f4e918b08b9/data-query/src/main/java/gov/va/api/health/dataquery/service/controller/patient/R4PatientTransformer.java#L399'

# --- Row 49: INSURED'S DOB ------------------------------------------------
$ws.Range('E49').Value = 'RelatedPerson'
$ws.Range('F49').Value = '.birthDate'

# --- Row 50: INSURED'S SSN ------------------------------------------------
$ws.Range('E50').Value = 'RelatedPerson'
$ws.Range('F50').Value = '.identifier'
$ws.Range('G50').Value = 'See R4PatientTransformer for identifier system/code details: https://github.com/department-of-veterans-affairs/health-apis-data-query/blob/a6b5bb26a88606753fdfff227cd23f4e918b08b9/data-query/src/main/java/gov/va/api/health/dataquery/service/controller/patient/R4PatientTransformer.java#L417'

# --- Row 54: INSURED'S SEX ------------------------------------------------
$ws.Range('E54').Value = 'RelatedPerson'
$ws.Range('F54').Value = '.extension'
$ws.Range('G54').Value = 'http://hl7.org/fhir/us/core/STU4/StructureDefinition-us-core-birthsex.html'

# --- Rows 71-77: SUBSCRIBER ADDRESS fields --------------------------------
$ws.Range('E71').Value = 'RelatedPerson'
$ws.Range('F71').Value = '.address'
$ws.Range('E72').Value = 'RelatedPerson'
$ws.Range('F72').Value = '.address'
$ws.Range('E73').Value = 'RelatedPerson'
$ws.Range('F73').Value = '.address'
$ws.Range('E74').Value = 'RelatedPerson'
$ws.Range('F74').Value = '.address'
$ws.Range('E75').Value = 'RelatedPerson'
$ws.Range('F75').Value = '.address'
$ws.Range('E76').Value = 'RelatedPerson'
$ws.Range('F76').Value = '.address'
$ws.Range('E77').Value = 'RelatedPerson'
$ws.Range('F77').Value = '.address'

# --- Row 78: SUBSCRIBER PHONE ---------------------------------------------
$ws.Range('E78').Value = 'RelatedPerson'
$ws.Range('F78').Value = '.telecom'

# --- Row 102: NAME OF INSURED ---------------------------------------------
$ws.Range('E102').Value = 'RelatedPerson'
$ws.Range('F102').Value = '.name.text'

# --- Widen column G so the long URLs/notes fit (best effort; the COM
#     width setter only resolves to the nearest 1/6 character) ------------
$ws.Columns('G').AutoFit()

# --- Move the view/selection to where the author ended up editing --------
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range('F103').Select()
